# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
# Both sheets share identical data, and the same set of rows changed in each.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3100
    5  = 2661
    7  = 138
    9  = 1398
    12 = 16
    13 = 1198
    14 = 4
    15 = 362
    16 = 327
    22 = 2559
    23 = 35
    24 = 292
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
